$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.403.03'
$ws.Range('E2').Value = '  +2.55%  '

$ws.Range('D3').Value = '2.426.89'
$ws.Range('E3').Value = '  +3.45%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').Value = "'556.72"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.39%  '

$ws.Range('D6').Value = "'144.31"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.64%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').Value = "'0.534"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.83%  '

$ws.Range('D9').Value = '2.428.18'
$ws.Range('E9').Value = '  +3.68%  '

$ws.Range('E10').Value = '  +5.51%  '

$ws.Range('E11').Value = '  -0.36%  '

$ws.Range('E12').Value = '  +2.25%  '

$ws.Range('D13').Value = "'0.351"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.76%  '

$ws.Range('D14').Value = "'26.38"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.16%  '

$ws.Range('E15').Value = '  +9.89%  '

$ws.Range('D16').Value = '2.864.63'
$ws.Range('E16').Value = '  +3.69%  '

$ws.Range('D17').Value = '62.309.03'
$ws.Range('E17').Value = '  +2.24%  '

$ws.Range('D18').Value = '2.427.37'
$ws.Range('E18').Value = '  +3.58%  '

$ws.Range('D19').Value = "'11.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.78%  '

$ws.Range('D20').Value = "'324.90"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.10%  '

$ws.Range('E21').Value = '  +1.70%  '

$ws.Range('E22').Value = '  +3.33%  '

$ws.Range('E23').Value = '  +0.40%  '

$ws.Range('D24').Value = "'1.80"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.51%  '

$ws.Range('D25').Value = "'64.99"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.78%  '

$ws.Range('D26').Value = "'9.10"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.82%  '

$ws.Range('D27').Value = "'573.24"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +15.29%  '

$ws.Range('D28').Value = '2.547.13'
$ws.Range('E28').Value = '  +3.85%  '

$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0949'
$ws.Range('E29').Value = '  +10.62%  '

$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.05%  '

$ws.Range('D31').Value = "'8.42"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.29%  '

$ws.Range('E33').Value = '  +2.34%  '

$ws.Range('D34').Value = "'1.86"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.40%  '

$ws.Range('D35').Value = "'1.58"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.15%  '

$ws.Range('D36').Value = "'5.73"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.31%  '

$ws.Range('D37').Value = "'4.86"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.98%  '

$ws.Range('E38').Value = '  +0.05%  '

$ws.Range('E39').Value = '  +2.70%  '

$ws.Range('E40').Value = '  +3.96%  '

$ws.Range('D41').Value = "'18.80"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.00%  '

$ws.Range('D42').Value = "'150.40"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.14%  '

$ws.Range('E43').Value = '  +0.06%  '

$ws.Range('E44').Value = '  +2.75%  '

$ws.Range('E45').Value = '  +15.63%  '

$ws.Range('D46').Value = "'151.31"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.12%  '

$ws.Range('E47').Value = '  +2.78%  '

$ws.Range('D48').Value = "'0.0544"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.17%  '

$ws.Range('D49').Value = "'20.48"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.72%  '

$ws.Range('E50').Value = '  +4.43%  '

$ws.Range('E51').Value = '  +2.28%  '
